# Add the new row of data (row 9), received from the bot, to the
# "Current_investments" worksheet.
#
# A9 is numeric (15). B9..F9 arrive as plain text (a date-looking string,
# a time-looking string, and plain digit strings) and must stay text,
# so the target range is pre-formatted as Text before the values are
# written -- otherwise Excel would silently reinterpret them as a date,
# a time and numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Current_investments")

$ws.Range("B9:F9").NumberFormat = "@"

$ws.Cells.Item(9, 1).Value = 15
$ws.Cells.Item(9, 2).Value = "2024.05.10"
$ws.Cells.Item(9, 3).Value = "23:00"
$ws.Cells.Item(9, 4).Value = "11"
$ws.Cells.Item(9, 5).Value = "22"
$ws.Cells.Item(9, 6).Value = "33"
